# chore: update Sheets via scheduled runner
# Refreshes cached market-price derived figures (currentAveragePrice*,
# LevePrice*/LeveProfit* columns H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# sheets to reflect newly polled data. Plain value overwrites - no formulas
# involved. A handful of rows gain or lose a trailing "LeveProfitHQ" (column N)
# cell depending on whether an HQ price is available this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 2149.5
$ws.Cells.Item(12, 9).Value = 1849.25
$ws.Cells.Item(12, 10).Value = 2449.75
$ws.Cells.Item(12, 11).Value = 1849.25
$ws.Cells.Item(12, 12).Value = 2449.75
$ws.Cells.Item(12, 13).Value = -1679.25
$ws.Cells.Item(12, 14).Value = -2789.75
$ws.Cells.Item(15, 8).Value = 498.6905
$ws.Cells.Item(15, 9).Value = 498.6905
$ws.Cells.Item(15, 11).Value = 1496.0715
$ws.Cells.Item(15, 13).Value = -1327.0715
$ws.Cells.Item(31, 8).Value = 93
$ws.Cells.Item(31, 9).Value = 93
$ws.Cells.Item(31, 11).Value = 279
$ws.Cells.Item(31, 13).Value = -49
$ws.Cells.Item(41, 8).Value = 411.2857
$ws.Cells.Item(41, 9).Value = 407.25
$ws.Cells.Item(41, 10).Value = 416.66666
$ws.Cells.Item(41, 11).Value = 407.25
$ws.Cells.Item(41, 12).Value = 416.66666
$ws.Cells.Item(41, 13).Value = 32.75
$ws.Cells.Item(41, 14).Value = -1296.66666
$ws.Cells.Item(76, 8).Value = 6588.3887
$ws.Cells.Item(76, 9).Value = 5762
$ws.Cells.Item(76, 11).Value = 5762
$ws.Cells.Item(76, 13).Value = -5447
$ws.Cells.Item(79, 8).Value = 6588.3887
$ws.Cells.Item(79, 9).Value = 5762
$ws.Cells.Item(79, 11).Value = 5762
$ws.Cells.Item(79, 13).Value = -4670
$ws.Cells.Item(106, 8).Value = 24260
$ws.Cells.Item(106, 9).Value = 27651.111
$ws.Cells.Item(106, 10).Value = 9000
$ws.Cells.Item(106, 11).Value = 27651.111
$ws.Cells.Item(106, 12).Value = 9000
$ws.Cells.Item(106, 13).Value = -27020.111
$ws.Cells.Item(106, 14).Value = -10262
$ws.Cells.Item(112, 8).Value = 2322.0557
$ws.Cells.Item(112, 10).Value = 2627.3572
$ws.Cells.Item(112, 12).Value = 7882.071599999999
$ws.Cells.Item(112, 14).Value = -10098.0716
$ws.Cells.Item(113, 8).Value = 2836.75
$ws.Cells.Item(113, 9).Value = 2582.3333
$ws.Cells.Item(113, 11).Value = 2582.3333
$ws.Cells.Item(113, 13).Value = 671.6667000000002
$ws.Cells.Item(135, 8).Value = 2499.8
$ws.Cells.Item(135, 10).Value = 3166.3333
$ws.Cells.Item(135, 12).Value = 28496.9997
$ws.Cells.Item(135, 14).Value = -33566.9997
$ws.Cells.Item(137, 8).Value = 2345.2068
$ws.Cells.Item(137, 9).Value = 1171.4615
$ws.Cells.Item(137, 11).Value = 3514.3845
$ws.Cells.Item(137, 13).Value = -964.3844999999997
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(110, 8).Value = 6399.7144
$ws.Cells.Item(110, 9).Value = 7899.8
$ws.Cells.Item(110, 10).Value = 2649.5
$ws.Cells.Item(110, 11).Value = 7899.8
$ws.Cells.Item(110, 12).Value = 2649.5
$ws.Cells.Item(110, 13).Value = -5854.8
$ws.Cells.Item(110, 14).Value = -6739.5
$ws.Cells.Item(122, 8).Value = 3371.8462
$ws.Cells.Item(122, 9).Value = 2047.7142
$ws.Cells.Item(122, 11).Value = 6143.142599999999
$ws.Cells.Item(122, 13).Value = -3693.142599999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1924
$ws.Cells.Item(20, 9).Value = 1887.25
$ws.Cells.Item(20, 10).Value = 1997.5
$ws.Cells.Item(20, 11).Value = 1887.25
$ws.Cells.Item(20, 12).Value = 1997.5
$ws.Cells.Item(20, 13).Value = -1640.25
$ws.Cells.Item(20, 14).Value = -2491.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 293.625
$ws.Cells.Item(22, 9).Value = 249.8
$ws.Cells.Item(22, 11).Value = 249.8
$ws.Cells.Item(22, 13).Value = 100.2
$ws.Cells.Item(134, 8).Value = 1588.8
$ws.Cells.Item(134, 9).Value = 1140.3235
$ws.Cells.Item(134, 10).Value = 2975
$ws.Cells.Item(134, 11).Value = 3420.9705
$ws.Cells.Item(134, 12).Value = 8925
$ws.Cells.Item(134, 13).Value = -885.9704999999999
$ws.Cells.Item(134, 14).Value = -13995
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1572704.9
$ws.Cells.Item(4, 9).Value = 2750362
$ws.Cells.Item(4, 11).Value = 8251086
$ws.Cells.Item(4, 13).Value = -8250974
$ws.Cells.Item(7, 8).Value = 5263316
$ws.Cells.Item(7, 9).Value = 8333396.5
$ws.Cells.Item(7, 10).Value = 320.7143
$ws.Cells.Item(7, 11).Value = 25000189.5
$ws.Cells.Item(7, 12).Value = 962.1428999999999
$ws.Cells.Item(7, 13).Value = -25000077.5
$ws.Cells.Item(7, 14).Value = -1186.1429
$ws.Cells.Item(12, 8).Value = 33.166668
$ws.Cells.Item(12, 10).Value = 33.166668
$ws.Cells.Item(12, 12).Value = 99.500004
$ws.Cells.Item(12, 14).Value = -445.500004
$ws.Cells.Item(50, 8).Value = 290.41666
$ws.Cells.Item(50, 9).Value = 174.28572
$ws.Cells.Item(50, 10).Value = 453
$ws.Cells.Item(50, 11).Value = 522.85716
$ws.Cells.Item(50, 12).Value = 1359
$ws.Cells.Item(50, 13).Value = -41.85716000000002
$ws.Cells.Item(50, 14).Value = -2321
$ws.Cells.Item(53, 8).Value = 290.41666
$ws.Cells.Item(53, 9).Value = 174.28572
$ws.Cells.Item(53, 10).Value = 453
$ws.Cells.Item(53, 11).Value = 522.85716
$ws.Cells.Item(53, 12).Value = 1359
$ws.Cells.Item(53, 13).Value = -41.85716000000002
$ws.Cells.Item(53, 14).Value = -2321
$ws.Cells.Item(55, 8).Value = 1740.25
$ws.Cells.Item(55, 10).Value = 1986
$ws.Cells.Item(55, 12).Value = 5958
$ws.Cells.Item(55, 14).Value = -6312
$ws.Cells.Item(59, 8).Value = 8000
$ws.Cells.Item(59, 10).Value = 8000
$ws.Cells.Item(59, 12).Value = 24000
$ws.Cells.Item(59, 14).Value = -25080
$ws.Cells.Item(62, 8).Value = 1533
$ws.Cells.Item(62, 9).Value = 1399.5
$ws.Cells.Item(62, 11).Value = 4198.5
$ws.Cells.Item(62, 13).Value = -3512.5
$ws.Cells.Item(65, 8).Value = 1533
$ws.Cells.Item(65, 9).Value = 1399.5
$ws.Cells.Item(65, 11).Value = 12595.5
$ws.Cells.Item(65, 13).Value = -9163.5
$ws.Cells.Item(109, 8).Value = 831.5
$ws.Cells.Item(109, 9).Value = 397.8
$ws.Cells.Item(109, 10).Value = 3000
$ws.Cells.Item(109, 11).Value = 1193.4
$ws.Cells.Item(109, 12).Value = 9000
$ws.Cells.Item(109, 13).Value = -153.4000000000001
$ws.Cells.Item(109, 14).Value = -11080
$ws.Cells.Item(126, 8).Value = 1950
$ws.Cells.Item(126, 9).Value = 1950
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 5850
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -910
$ws.Cells.Item(126, 14).Value = $null
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 281.46155
$ws.Cells.Item(2, 9).Value = 76.625
$ws.Cells.Item(2, 11).Value = 76.625
$ws.Cells.Item(2, 13).Value = 36.375
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 12).Value = 0
$ws.Cells.Item(86, 14).Value = $null
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 10).Value = 0
$ws.Cells.Item(89, 12).Value = 0
$ws.Cells.Item(89, 14).Value = $null
$ws.Cells.Item(122, 8).Value = 480615.3
$ws.Cells.Item(122, 9).Value = 57231.445
$ws.Cells.Item(122, 11).Value = 171694.335
$ws.Cells.Item(122, 13).Value = -169244.335
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 360.3913
$ws.Cells.Item(55, 10).Value = 487.6
$ws.Cells.Item(55, 12).Value = 487.6
$ws.Cells.Item(55, 14).Value = -833.6
$ws.Cells.Item(61, 8).Value = 6060.875
$ws.Cells.Item(61, 9).Value = 6355.143
$ws.Cells.Item(61, 11).Value = 6355.143
$ws.Cells.Item(61, 13).Value = -6153.143
$ws.Cells.Item(82, 8).Value = 3455
$ws.Cells.Item(82, 9).Value = 3437
$ws.Cells.Item(82, 11).Value = 3437
$ws.Cells.Item(82, 13).Value = -3076
$ws.Cells.Item(85, 8).Value = 3455
$ws.Cells.Item(85, 9).Value = 3437
$ws.Cells.Item(85, 11).Value = 3437
$ws.Cells.Item(85, 13).Value = -2189
$ws.Cells.Item(113, 8).Value = 6060.875
$ws.Cells.Item(113, 9).Value = 6355.143
$ws.Cells.Item(113, 11).Value = 6355.143
$ws.Cells.Item(113, 13).Value = -4185.143
$ws.Cells.Item(122, 8).Value = 3446
$ws.Cells.Item(122, 9).Value = 3446
$ws.Cells.Item(122, 11).Value = 10338
$ws.Cells.Item(122, 13).Value = -7888
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 1064.4736
$ws.Cells.Item(81, 9).Value = 979.1667
$ws.Cells.Item(81, 10).Value = 2600
$ws.Cells.Item(81, 11).Value = 1958.3334
$ws.Cells.Item(81, 12).Value = 5200
$ws.Cells.Item(81, 13).Value = -897.3334
$ws.Cells.Item(81, 14).Value = -7322
$ws.Cells.Item(84, 8).Value = 1064.4736
$ws.Cells.Item(84, 9).Value = 979.1667
$ws.Cells.Item(84, 10).Value = 2600
$ws.Cells.Item(84, 11).Value = 9791.666999999999
$ws.Cells.Item(84, 12).Value = 26000
$ws.Cells.Item(84, 13).Value = -4487.666999999999
$ws.Cells.Item(84, 14).Value = -36608
$ws.Cells.Item(122, 8).Value = 1347.1
$ws.Cells.Item(122, 9).Value = 1286.421
$ws.Cells.Item(122, 11).Value = 3859.263
$ws.Cells.Item(122, 13).Value = -1409.263
$ws.Cells.Item(124, 8).Value = 29199
$ws.Cells.Item(124, 10).Value = 29199
$ws.Cells.Item(124, 12).Value = 29199
$ws.Cells.Item(124, 14).Value = -39019
$ws.Cells.Item(132, 8).Value = 1871.258
$ws.Cells.Item(132, 9).Value = 1438
$ws.Cells.Item(132, 10).Value = 3356.7144
$ws.Cells.Item(132, 11).Value = 4314
$ws.Cells.Item(132, 12).Value = 10070.1432
$ws.Cells.Item(132, 13).Value = -1784
$ws.Cells.Item(132, 14).Value = -15130.1432
